$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates pulled from the latest coinranking.com snapshot.
# Columns B/C/D hold text (coin name / link / price-as-text), so we
# force Text number format before writing to stop Excel from
# reinterpreting values like "582.04" or "0.0000161" as numbers,
# then restore the Normal style so no formatting residue is left behind.
$updates = @(
    @{ Cell = 'D2'; Value = '65.650.17'; ForceText = $True }
    @{ Cell = 'E2'; Value = '  -1.38%  '; ForceText = $False }
    @{ Cell = 'D3'; Value = '3.269.79'; ForceText = $True }
    @{ Cell = 'E3'; Value = '  -1.12%  '; ForceText = $False }
    @{ Cell = 'E4'; Value = '  -0.03%  '; ForceText = $False }
    @{ Cell = 'D5'; Value = '582.04'; ForceText = $True }
    @{ Cell = 'E5'; Value = '  +1.81%  '; ForceText = $False }
    @{ Cell = 'D6'; Value = '178.48'; ForceText = $True }
    @{ Cell = 'E6'; Value = '  -2.12%  '; ForceText = $False }
    @{ Cell = 'D7'; Value = '0.639'; ForceText = $True }
    @{ Cell = 'E7'; Value = '  +6.89%  '; ForceText = $False }
    @{ Cell = 'E8'; Value = '  +0.01%  '; ForceText = $False }
    @{ Cell = 'E9'; Value = '  -3.82%  '; ForceText = $False }
    @{ Cell = 'D10'; Value = '6.72'; ForceText = $True }
    @{ Cell = 'E10'; Value = '  +1.43%  '; ForceText = $False }
    @{ Cell = 'D11'; Value = '0.400'; ForceText = $True }
    @{ Cell = 'E11'; Value = '  -0.74%  '; ForceText = $False }
    @{ Cell = 'D12'; Value = '3.835.88'; ForceText = $True }
    @{ Cell = 'E12'; Value = '  -1.17%  '; ForceText = $False }
    @{ Cell = 'E13'; Value = '  -4.50%  '; ForceText = $False }
    @{ Cell = 'D14'; Value = '65.762.99'; ForceText = $True }
    @{ Cell = 'E14'; Value = '  -1.34%  '; ForceText = $False }
    @{ Cell = 'D15'; Value = '25.86'; ForceText = $True }
    @{ Cell = 'E15'; Value = '  -4.69%  '; ForceText = $False }
    @{ Cell = 'B16'; Value = 'WrappedEther'; ForceText = $True }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; ForceText = $True }
    @{ Cell = 'D16'; Value = '3.272.22'; ForceText = $True }
    @{ Cell = 'E16'; Value = '  -1.09%  '; ForceText = $False }
    @{ Cell = 'B17'; Value = 'ShibaInu'; ForceText = $True }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; ForceText = $True }
    @{ Cell = 'D17'; Value = '0.0000161'; ForceText = $True }
    @{ Cell = 'E17'; Value = '  -3.38%  '; ForceText = $False }
    @{ Cell = 'D18'; Value = '424.19'; ForceText = $True }
    @{ Cell = 'E18'; Value = '  -1.77%  '; ForceText = $False }
    @{ Cell = 'D19'; Value = '13.11'; ForceText = $True }
    @{ Cell = 'E19'; Value = '  -4.03%  '; ForceText = $False }
    @{ Cell = 'D20'; Value = '5.46'; ForceText = $True }
    @{ Cell = 'E20'; Value = '  -3.54%  '; ForceText = $False }
    @{ Cell = 'D21'; Value = '7.32'; ForceText = $True }
    @{ Cell = 'E21'; Value = '  -3.76%  '; ForceText = $False }
    @{ Cell = 'B22'; Value = 'Dai'; ForceText = $True }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; ForceText = $True }
    @{ Cell = 'D22'; Value = '1.00'; ForceText = $True }
    @{ Cell = 'E22'; Value = '  +0.29%  '; ForceText = $False }
    @{ Cell = 'B23'; Value = 'Litecoin'; ForceText = $True }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; ForceText = $True }
    @{ Cell = 'D23'; Value = '71.47'; ForceText = $True }
    @{ Cell = 'E23'; Value = '  -2.98%  '; ForceText = $False }
    @{ Cell = 'E24'; Value = '  +0.27%  '; ForceText = $False }
    @{ Cell = 'D25'; Value = '3.421.73'; ForceText = $True }
    @{ Cell = 'E25'; Value = '  -0.99%  '; ForceText = $False }
    @{ Cell = 'D26'; Value = '0.506'; ForceText = $True }
    @{ Cell = 'E26'; Value = '  -1.81%  '; ForceText = $False }
    @{ Cell = 'E27'; Value = '  +1.40%  '; ForceText = $False }
    @{ Cell = 'D28'; Value = '0.0000111'; ForceText = $True }
    @{ Cell = 'E28'; Value = '  -5.46%  '; ForceText = $False }
    @{ Cell = 'D29'; Value = '8.80'; ForceText = $True }
    @{ Cell = 'E29'; Value = '  -2.47%  '; ForceText = $False }
    @{ Cell = 'E30'; Value = '  -0.02%  '; ForceText = $False }
    @{ Cell = 'E31'; Value = '  -0.81%  '; ForceText = $False }
    @{ Cell = 'D32'; Value = '22.13'; ForceText = $True }
    @{ Cell = 'E32'; Value = '  -2.67%  '; ForceText = $False }
    @{ Cell = 'D33'; Value = '0.999'; ForceText = $True }
    @{ Cell = 'E33'; Value = '  +0.04%  '; ForceText = $False }
    @{ Cell = 'D34'; Value = '5.11'; ForceText = $True }
    @{ Cell = 'E34'; Value = '  -4.00%  '; ForceText = $False }
    @{ Cell = 'D35'; Value = '6.54'; ForceText = $True }
    @{ Cell = 'E35'; Value = '  -3.43%  '; ForceText = $False }
    @{ Cell = 'D37'; Value = '158.97'; ForceText = $True }
    @{ Cell = 'E38'; Value = '  -5.89%  '; ForceText = $False }
    @{ Cell = 'D39'; Value = '1.78'; ForceText = $True }
    @{ Cell = 'E39'; Value = '  -3.19%  '; ForceText = $False }
    @{ Cell = 'D40'; Value = '26.18'; ForceText = $True }
    @{ Cell = 'E40'; Value = '  -3.56%  '; ForceText = $False }
    @{ Cell = 'D41'; Value = '2.779.75'; ForceText = $True }
    @{ Cell = 'E41'; Value = '  -0.78%  '; ForceText = $False }
    @{ Cell = 'D42'; Value = '0.763'; ForceText = $True }
    @{ Cell = 'E42'; Value = '  -3.29%  '; ForceText = $False }
    @{ Cell = 'D43'; Value = '4.29'; ForceText = $True }
    @{ Cell = 'E43'; Value = '  -3.28%  '; ForceText = $False }
    @{ Cell = 'D44'; Value = '39.81'; ForceText = $True }
    @{ Cell = 'E44'; Value = '  -0.86%  '; ForceText = $False }
    @{ Cell = 'D45'; Value = '0.0653'; ForceText = $True }
    @{ Cell = 'E45'; Value = '  -3.15%  '; ForceText = $False }
    @{ Cell = 'D46'; Value = '5.83'; ForceText = $True }
    @{ Cell = 'E46'; Value = '  -5.65%  '; ForceText = $False }
    @{ Cell = 'D47'; Value = '2.26'; ForceText = $True }
    @{ Cell = 'E47'; Value = '  -3.59%  '; ForceText = $False }
    @{ Cell = 'D48'; Value = '313.33'; ForceText = $True }
    @{ Cell = 'E48'; Value = '  -2.08%  '; ForceText = $False }
    @{ Cell = 'D49'; Value = '22.95'; ForceText = $True }
    @{ Cell = 'E49'; Value = '  -5.70%  '; ForceText = $False }
    @{ Cell = 'D50'; Value = '0.0265'; ForceText = $True }
    @{ Cell = 'E50'; Value = '  -2.24%  '; ForceText = $False }
    @{ Cell = 'D51'; Value = '0.103'; ForceText = $True }
    @{ Cell = 'E51'; Value = '  +3.16%  '; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
